$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ASV_rank (column G) from 41 to 42 for rows 9-12
$ws.Range("G9:G12").Value = 42
